$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '57.207.08'
Set-TextValue 'E2' '  +4.29%  '

# Row 3
Set-TextValue 'D3' '2.511.16'
Set-TextValue 'E3' '  +2.30%  '

# Row 4
Set-TextValue 'E4' '  -0.15%  '

# Row 5
Set-TextValue 'D5' '496.20'
Set-TextValue 'E5' '  +3.79%  '

# Row 6
Set-TextValue 'D6' '153.98'
Set-TextValue 'E6' '  +11.15%  '

# Row 7
Set-TextValue 'D7' '0.996'
Set-TextValue 'E7' '  -0.34%  '

# Row 8
Set-TextValue 'D8' '0.517'
Set-TextValue 'E8' '  +3.56%  '

# Row 9
Set-TextValue 'D9' '2.525.86'
Set-TextValue 'E9' '  +3.09%  '

# Row 10
Set-TextValue 'D10' '0.101'
Set-TextValue 'E10' '  +5.54%  '

# Row 11
Set-TextValue 'E11' '  +5.77%  '

# Row 12
Set-TextValue 'D12' '0.339'
Set-TextValue 'E12' '  +4.80%  '

# Row 13
Set-TextValue 'E13' '  +1.57%  '

# Row 14
Set-TextValue 'D14' '2.949.89'
Set-TextValue 'E14' '  +2.83%  '

# Row 15
Set-TextValue 'D15' '57.339.85'
Set-TextValue 'E15' '  +4.12%  '

# Row 16
Set-TextValue 'D16' '21.38'
Set-TextValue 'E16' '  +5.03%  '

# Row 17
Set-TextValue 'D17' '0.0000139'
Set-TextValue 'E17' '  +3.55%  '

# Row 18
Set-TextValue 'D18' '2.517.50'
Set-TextValue 'E18' '  +2.91%  '

# Row 19
Set-TextValue 'D19' '4.60'
Set-TextValue 'E19' '  +6.01%  '

# Row 20
Set-TextValue 'D20' '10.34'
Set-TextValue 'E20' '  +4.53%  '

# Row 21
Set-TextValue 'D21' '324.41'
Set-TextValue 'E21' '  +3.76%  '

# Row 22
Set-TextValue 'D22' '0.999'
Set-TextValue 'E22' '  +0.35%  '

# Row 23
Set-TextValue 'D23' '5.95'
Set-TextValue 'E23' '  +6.07%  '

# Row 24
Set-TextValue 'D24' '58.50'
Set-TextValue 'E24' '  +2.36%  '

# Row 25
Set-TextValue 'D25' '0.411'
Set-TextValue 'E25' '  +2.12%  '

# Row 26
Set-TextValue 'B26' 'Kaspa'
Set-TextValue 'C26' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D26' '0.166'
Set-TextValue 'E26' '  +3.11%  '

# Row 27
Set-TextValue 'B27' 'Binance-PegBSC-USD'
Set-TextValue 'C27' 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue 'D27' '1.00'
Set-TextValue 'E27' '  -0.12%  '

# Row 28
Set-TextValue 'D28' '2.617.70'
Set-TextValue 'E28' '  +2.70%  '

# Row 29
Set-TextValue 'D29' '7.64'
Set-TextValue 'E29' '  +4.65%  '

# Row 30
Set-TextValue 'D30' '0.0₃0832'
Set-TextValue 'E30' '  +8.41%  '

# Row 31
Set-TextValue 'E31' '  -0.21%  '

# Row 32
Set-TextValue 'D32' '151.59'
Set-TextValue 'E32' '  +2.00%  '

# Row 33
Set-TextValue 'B33' 'EthereumClassic'
Set-TextValue 'C33' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D33' '18.42'
Set-TextValue 'E33' '  +3.16%  '

# Row 34
Set-TextValue 'B34' 'PancakeSwap'
Set-TextValue 'C34' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D34' '1.54'
Set-TextValue 'E34' '  +4.88%  '

# Row 35
Set-TextValue 'D35' '5.29'
Set-TextValue 'E35' '  +2.96%  '

# Row 36
Set-TextValue 'D36' '3.84'
Set-TextValue 'E36' '  +7.28%  '

# Row 37
Set-TextValue 'B37' 'Fetch.AI'
Set-TextValue 'C37' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D37' '0.903'
Set-TextValue 'E37' '  +6.30%  '

# Row 38
Set-TextValue 'B38' 'ImmutableX'
Set-TextValue 'C38' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D38' '1.16'
Set-TextValue 'E38' '  +4.83%  '

# Row 39
Set-TextValue 'D39' '1.42'
Set-TextValue 'E39' '  +10.64%  '

# Row 40
Set-TextValue 'D40' '34.43'
Set-TextValue 'E40' '  +3.97%  '

# Row 41
Set-TextValue 'D41' '3.54'
Set-TextValue 'E41' '  +4.75%  '

# Row 42
Set-TextValue 'D42' '0.620'
Set-TextValue 'E42' '  +3.39%  '

# Row 43
Set-TextValue 'D43' '0.0564'
Set-TextValue 'E43' '  +4.18%  '

# Row 44
Set-TextValue 'D44' '0.994'
Set-TextValue 'E44' '  -0.20%  '

# Row 45
Set-TextValue 'D45' '4.93'
Set-TextValue 'E45' '  +7.54%  '

# Row 46
Set-TextValue 'D46' '269.13'
Set-TextValue 'E46' '  +6.08%  '

# Row 47
Set-TextValue 'D47' '0.0949'
Set-TextValue 'E47' '  +6.49%  '

# Row 48
Set-TextValue 'D48' '0.0231'
Set-TextValue 'E48' '  +4.48%  '

# Row 49
Set-TextValue 'D49' '10.23'
Set-TextValue 'E49' '  +0.74%  '

# Row 50
Set-TextValue 'D50' '18.11'
Set-TextValue 'E50' '  +6.86%  '

# Row 51
Set-TextValue 'D51' '1.899.52'
Set-TextValue 'E51' '  -1.21%  '
